$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number for every data row
# (rows 2-351), all currently set to 45178 (2023-09-09). This mirrors the
# automatic daily-refresh commit: bump each of those cells to 45179
# (2023-09-10), one day later, leaving every other column untouched.
$range = $ws.Range("C2:C351")
for ($i = 1; $i -le $range.Rows.Count; $i++) {
    $cell = $range.Cells.Item($i, 1)
    if ($cell.Value2 -eq 45178) {
        $cell.Value2 = 45179
    }
}
